$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")

$fv2310Headers = @("Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310")
$fv2404Headers = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}
# column 11 (K) is "diff" - stays the same
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# 2) Turn the used range into an Excel Table ("Table1")
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U69"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3) Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
